$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -6
$ws.Range("F9").Value = -4
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = -4
